# Add a new worksheet "Торец" (Edge/Side) by duplicating the loading-side
# material price table that already lives on "Задник" (rows 1:41), placing
# the new sheet right after "Каркас", then trimming it down to just that
# table.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Задник")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate "Задник" as a new sheet positioned after the last sheet ("Каркас").
$src.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Торец"

# Only rows 1:41 (the loading-side material table) are needed; drop the rest.
$newSheet.Rows("42:71").Delete()

$newSheet.Columns("B:B").ColumnWidth = 8.17

$newSheet.Range("A1:B41").Select()
$newSheet.Cells.Item(12, 4).Select()
